$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.448696851730347
$ws.Range("B1").Value = 1.517573475837708
$ws.Range("C1").Value = 1.622479915618896
$ws.Range("D1").Value = 2.282655000686646
$ws.Range("E1").Value = 3.755235195159912
